$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 33
$wsALC.Range("H33").Value = 323.19232
$wsALC.Range("I33").Value = 325.125
$wsALC.Range("J33").Value = 300
$wsALC.Range("K33").Value = 325.125
$wsALC.Range("L33").Value = 300
$wsALC.Range("M33").Value = -96.125
$wsALC.Range("N33").Value = -758

# ALC row 129
$wsALC.Range("H129").Value = 1820.6154
$wsALC.Range("J129").Value = 2051.7856
$wsALC.Range("L129").Value = 6155.3568
$wsALC.Range("N129").Value = -16155.3568

# ALC row 135
$wsALC.Range("H135").Value = 752.625
$wsALC.Range("I135").Value = 781.9091
$wsALC.Range("J135").Value = 430.5
$wsALC.Range("K135").Value = 7037.1819
$wsALC.Range("L135").Value = 3874.5
$wsALC.Range("M135").Value = -4502.1819
$wsALC.Range("N135").Value = -8944.5

# ALC row 137
$wsALC.Range("H137").Value = 5263937
$wsALC.Range("I137").Value = 582.3
$wsALC.Range("J137").Value = 11112109
$wsALC.Range("K137").Value = 1746.9
$wsALC.Range("L137").Value = 33336327
$wsALC.Range("M137").Value = 803.1000000000001
$wsALC.Range("N137").Value = -33341427

# ALC row 138
$wsALC.Range("H138").Value = 6668615
$wsALC.Range("I138").Value = 12821815
$wsALC.Range("J138").Value = 2647.9167
$wsALC.Range("K138").Value = 38465445
$wsALC.Range("L138").Value = 7943.750100000001
$wsALC.Range("M138").Value = -38460305
$wsALC.Range("N138").Value = -18223.7501

# ALC row 141
$wsALC.Range("H141").Value = 1275.0256
$wsALC.Range("I141").Value = 1135.0286
$wsALC.Range("J141").Value = 2500
$wsALC.Range("K141").Value = 3405.0858
$wsALC.Range("L141").Value = 7500
$wsALC.Range("M141").Value = 1774.9142
$wsALC.Range("N141").Value = -17860

# ARM row 32
$wsARM.Range("H32").Value = 8803.710999999999
$wsARM.Range("I32").Value = 8495.066000000001
$wsARM.Range("J32").Value = 10346.934
$wsARM.Range("K32").Value = 8495.066000000001
$wsARM.Range("L32").Value = 10346.934
$wsARM.Range("M32").Value = -8208.066000000001
$wsARM.Range("N32").Value = -10920.934

# ARM row 61
$wsARM.Range("H61").Value = 11365386
$wsARM.Range("I61").Value = 13890556
$wsARM.Range("J61").Value = 2118.75
$wsARM.Range("K61").Value = 13890556
$wsARM.Range("L61").Value = 2118.75
$wsARM.Range("M61").Value = -13890344
$wsARM.Range("N61").Value = -2542.75

# ARM row 74
$wsARM.Range("H74").Value = 12503104
$wsARM.Range("I74").Value = 16668755
$wsARM.Range("J74").Value = 6152.8
$wsARM.Range("K74").Value = 16668755
$wsARM.Range("L74").Value = 6152.8
$wsARM.Range("M74").Value = -16667881
$wsARM.Range("N74").Value = -7900.8

# ARM row 77
$wsARM.Range("H77").Value = 12503104
$wsARM.Range("I77").Value = 16668755
$wsARM.Range("J77").Value = 6152.8
$wsARM.Range("K77").Value = 83343775
$wsARM.Range("L77").Value = 30764
$wsARM.Range("M77").Value = -83339407
$wsARM.Range("N77").Value = -39500

# ARM row 132
$wsARM.Range("H132").Value = 6412783.5
$wsARM.Range("I132").Value = 8335430.5
$wsARM.Range("J132").Value = 3960
$wsARM.Range("K132").Value = 25006291.5
$wsARM.Range("L132").Value = 11880
$wsARM.Range("M132").Value = -25003761.5
$wsARM.Range("N132").Value = -16940

# ARM row 136
$wsARM.Range("H136").Value = 11365386
$wsARM.Range("I136").Value = 13890556
$wsARM.Range("J136").Value = 2118.75
$wsARM.Range("K136").Value = 41671668
$wsARM.Range("L136").Value = 6356.25
$wsARM.Range("M136").Value = -41669118
$wsARM.Range("N136").Value = -11456.25

# BSM row 134
$wsBSM.Range("H134").Value = 3260.5952
$wsBSM.Range("I134").Value = 1930.5186
$wsBSM.Range("J134").Value = 5654.7334
$wsBSM.Range("K134").Value = 5791.5558
$wsBSM.Range("L134").Value = 16964.2002
$wsBSM.Range("M134").Value = -3256.5558
$wsBSM.Range("N134").Value = -22034.2002

# CRP row 31
$wsCRP.Range("H31").Value = 10759082
$wsCRP.Range("I31").Value = 7056.826
$wsCRP.Range("J31").Value = 41671150
$wsCRP.Range("K31").Value = 7056.826
$wsCRP.Range("L31").Value = 41671150
$wsCRP.Range("M31").Value = -6761.826
$wsCRP.Range("N31").Value = -41671740

# CRP row 34
$wsCRP.Range("H34").Value = 10759082
$wsCRP.Range("I34").Value = 7056.826
$wsCRP.Range("J34").Value = 41671150
$wsCRP.Range("K34").Value = 7056.826
$wsCRP.Range("L34").Value = 41671150
$wsCRP.Range("M34").Value = -6854.826
$wsCRP.Range("N34").Value = -41671554

# CRP row 132
$wsCRP.Range("H132").Value = 10418619
$wsCRP.Range("I132").Value = 13515248
$wsCRP.Range("J132").Value = 2684.5454
$wsCRP.Range("K132").Value = 40545744
$wsCRP.Range("L132").Value = 8053.6362
$wsCRP.Range("M132").Value = -40543214
$wsCRP.Range("N132").Value = -13113.6362

# CUL row 92
$wsCUL.Range("H92").Value = 850.75
$wsCUL.Range("I92").Value = 0
$wsCUL.Range("J92").Value = 850.75
$wsCUL.Range("K92").Value = 0
$wsCUL.Range("L92").Value = 2552.25
$wsCUL.Range("M92").ClearContents()
$wsCUL.Range("N92").Value = -5048.25

# CUL row 131
$wsCUL.Range("H131").Value = 815.09
$wsCUL.Range("I131").Value = 465.45456
$wsCUL.Range("J131").Value = 858.30334
$wsCUL.Range("K131").Value = 1396.36368
$wsCUL.Range("L131").Value = 2574.91002
$wsCUL.Range("M131").Value = 3643.63632
$wsCUL.Range("N131").Value = -12654.91002

# GSM row 113
$wsGSM.Range("H113").Value = 22392.064
$wsGSM.Range("I113").Value = 28455.861
$wsGSM.Range("J113").Value = 562.4
$wsGSM.Range("K113").Value = 28455.861
$wsGSM.Range("L113").Value = 562.4
$wsGSM.Range("M113").Value = -26285.861
$wsGSM.Range("N113").Value = -4902.4

# GSM row 126
$wsGSM.Range("H126").Value = 4440
$wsGSM.Range("I126").Value = 2585
$wsGSM.Range("J126").Value = 6913.3335
$wsGSM.Range("K126").Value = 7755
$wsGSM.Range("L126").Value = 20740.0005
$wsGSM.Range("M126").Value = -5285
$wsGSM.Range("N126").Value = -25680.0005

# LTW row 7
$wsLTW.Range("H7").Value = 7490.091
$wsLTW.Range("I7").Value = 10560.25
$wsLTW.Range("J7").Value = 5735.7144
$wsLTW.Range("K7").Value = 10560.25
$wsLTW.Range("L7").Value = 5735.7144
$wsLTW.Range("M7").Value = -10448.25
$wsLTW.Range("N7").Value = -5959.7144

# LTW row 55
$wsLTW.Range("H55").Value = 555.4545000000001
$wsLTW.Range("I55").Value = 540
$wsLTW.Range("J55").Value = 568.3333
$wsLTW.Range("K55").Value = 540
$wsLTW.Range("L55").Value = 568.3333
$wsLTW.Range("M55").Value = -367
$wsLTW.Range("N55").Value = -914.3333

# LTW row 126
$wsLTW.Range("H126").Value = 7490.091
$wsLTW.Range("I126").Value = 10560.25
$wsLTW.Range("J126").Value = 5735.7144
$wsLTW.Range("K126").Value = 31680.75
$wsLTW.Range("L126").Value = 17207.1432
$wsLTW.Range("M126").Value = -29210.75
$wsLTW.Range("N126").Value = -22147.1432

# LTW row 136
$wsLTW.Range("H136").Value = 13162138
$wsLTW.Range("I136").Value = 17858598
$wsLTW.Range("J136").Value = 12050.5
$wsLTW.Range("K136").Value = 53575794
$wsLTW.Range("L136").Value = 36151.5
$wsLTW.Range("M136").Value = -53573244
$wsLTW.Range("N136").Value = -41251.5

# LTW row 138
$wsLTW.Range("H138").Value = 88876.86
$wsLTW.Range("J138").Value = 88876.86
$wsLTW.Range("L138").Value = 88876.86
$wsLTW.Range("N138").Value = -99156.86

# WVR row 126
$wsWVR.Range("H126").Value = 3641.1875
$wsWVR.Range("I126").Value = 2573.0908
$wsWVR.Range("J126").Value = 5991
$wsWVR.Range("K126").Value = 7719.2724
$wsWVR.Range("L126").Value = 17973
$wsWVR.Range("M126").Value = -5249.2724
$wsWVR.Range("N126").Value = -22913

# WVR row 132
$wsWVR.Range("H132").Value = 1479.0819
$wsWVR.Range("I132").Value = 1232.7
$wsWVR.Range("J132").Value = 2599
$wsWVR.Range("K132").Value = 3698.1
$wsWVR.Range("L132").Value = 7797
$wsWVR.Range("M132").Value = -1168.1
$wsWVR.Range("N132").Value = -12857

# WVR row 136
$wsWVR.Range("H136").Value = 843.0714
$wsWVR.Range("I136").Value = 774.6667
$wsWVR.Range("J136").Value = 966.2
$wsWVR.Range("K136").Value = 2324.0001
$wsWVR.Range("L136").Value = 2898.6
$wsWVR.Range("M136").Value = 225.9998999999998
$wsWVR.Range("N136").Value = -7998.6

